$d = $word.ActiveDocument

# The team roster table's 5th data row (the first of the two trailing
# blank rows) is currently empty across all three cells. Fill in the
# new attendee's roll number, name, and position title, matching the
# same run/paragraph formatting (theme fonts + text color) used by the
# other populated rows in the table.

$table = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'
$rPr = '<w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:color w:val="2D3B45"/></w:rPr>'
$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="0"/><w:outlineLvl w:val="0"/>' + $rPr + '</w:pPr>'

function Fill-Cell($rowIndex, $cellIndex, $paraId, $textId, $text) {
    # Table.Cell(row, col) addresses the cell directly (unlike
    # Rows.Item(row).Cells.Item(col), whose .Range can cause InsertXML
    # to insert a sibling paragraph instead of replacing the existing
    # one for any cell after the first in a row).
    $cellRange = $table.Cell($rowIndex, $cellIndex).Range
    $attrs = 'w14:paraId="' + $paraId + '" w14:textId="' + $textId + '" w:rsidR="00685ECC" w:rsidRPr="00A6112F" w:rsidRDefault="00685ECC" w:rsidP="00685ECC"'
    $run = '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    $xml = '<w:p ' + $wNs + ' ' + $attrs + '>' + $pPr + $run + '</w:p>'
    $cellRange.InsertXML($xml)
}

Fill-Cell 5 1 "142B03E6" "61F2E110" "19"
Fill-Cell 5 2 "51413AD6" "3DC5A25F" "Hunter Malinowski"
Fill-Cell 5 3 "4E7B0AFE" "0AECE97A" "Developer"
